$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the existing header cells
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I0 (I) and IF (J) data columns for rows 2-38
$data = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 3, 4),
    @(5, 6, 6),
    @(6, 7, 8),
    @(7, 7, 7),
    @(8, 7, 8),
    @(9, 5, 6),
    @(10, 3, 4),
    @(11, 6, 6),
    @(12, 8, 9),
    @(13, 7, 8),
    @(14, 6, 6),
    @(15, 7, 7),
    @(16, 8, 8),
    @(17, 8, 8),
    @(18, 8, 9),
    @(19, 8, 9),
    @(20, 5, 5),
    @(21, 3, 4),
    @(22, 8, 8),
    @(23, 5, 6),
    @(24, 4, 5),
    @(25, 9, 9),
    @(26, 6, 6),
    @(27, 9, 9),
    @(28, 7, 7),
    @(29, 8, 8),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 5, 5),
    @(33, 9, 9),
    @(34, 6, 6),
    @(35, 8, 8),
    @(36, 5, 6),
    @(37, 5, 5),
    @(38, 3, 3)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
